# Auto-generated Excel COM-interop script
# Updates market-price columns (H-N) on 28 rows across 8 sheets
# per the "chore: update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2020.4445
$ws.Range("I70").Value = 1424.5
$ws.Range("K70").Value = 4273.5
$ws.Range("M70").Value = -4003.5
# Row 73
$ws.Range("H73").Value = 2020.4445
$ws.Range("I73").Value = 1424.5
$ws.Range("K73").Value = 4273.5
$ws.Range("M73").Value = -3337.5
# Row 98
$ws.Range("H98").Value = 2946.12
$ws.Range("I98").Value = 3193.3914
$ws.Range("J98").Value = 102.5
$ws.Range("K98").Value = 3193.3914
$ws.Range("L98").Value = 102.5
$ws.Range("M98").Value = -1695.3914
$ws.Range("N98").Value = -3098.5
# Row 122
$ws.Range("H122").Value = 2946.12
$ws.Range("I122").Value = 3193.3914
$ws.Range("J122").Value = 102.5
$ws.Range("K122").Value = 9580.174199999999
$ws.Range("L122").Value = 307.5
$ws.Range("M122").Value = -7130.174199999999
$ws.Range("N122").Value = -5207.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3418.4717
$ws.Range("I32").Value = 3076.7917
$ws.Range("J32").Value = 6698.6
$ws.Range("K32").Value = 3076.7917
$ws.Range("L32").Value = 6698.6
$ws.Range("M32").Value = -2789.7917
$ws.Range("N32").Value = -7272.6

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3083.5
$ws.Range("I86").Value = 3154.2083
$ws.Range("J86").Value = 2942.0833
$ws.Range("K86").Value = 3154.2083
$ws.Range("L86").Value = 2942.0833
$ws.Range("M86").Value = -2031.2083
$ws.Range("N86").Value = -5188.0833
# Row 89
$ws.Range("H89").Value = 3083.5
$ws.Range("I89").Value = 3154.2083
$ws.Range("J89").Value = 2942.0833
$ws.Range("K89").Value = 15771.0415
$ws.Range("L89").Value = 14710.4165
$ws.Range("M89").Value = -10155.0415
$ws.Range("N89").Value = -25942.4165
# Row 94
$ws.Range("H94").Value = 8929047
$ws.Range("I94").Value = 12500432
$ws.Range("K94").Value = 12500432
$ws.Range("M94").Value = -12499981
# Row 105
$ws.Range("H105").Value = 100002440
$ws.Range("I105").Value = 100002440
$ws.Range("K105").Value = 100002440
$ws.Range("M105").Value = -100000693
# Row 134
$ws.Range("H134").Value = 6852.909
$ws.Range("I134").Value = 1160.4375
$ws.Range("K134").Value = 3481.3125
$ws.Range("M134").Value = -946.3125

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 20002750
$ws.Range("J62").Value = 66668500
$ws.Range("L62").Value = 66668500
$ws.Range("N62").Value = -66669748
# Row 65
$ws.Range("H65").Value = 20002750
$ws.Range("J65").Value = 66668500
$ws.Range("L65").Value = 333342500
$ws.Range("N65").Value = -333348740
# Row 99
$ws.Range("H99").Value = 1876.5714
$ws.Range("I99").Value = 1870.3334
$ws.Range("J99").Value = 1914
$ws.Range("K99").Value = 1870.3334
$ws.Range("L99").Value = 1914
$ws.Range("M99").Value = -372.3334
$ws.Range("N99").Value = -4910
# Row 126
$ws.Range("H126").Value = 1876.5714
$ws.Range("I126").Value = 1870.3334
$ws.Range("J126").Value = 1914
$ws.Range("K126").Value = 5611.0002
$ws.Range("L126").Value = 5742
$ws.Range("M126").Value = -3141.0002
$ws.Range("N126").Value = -10682
# Row 132
$ws.Range("H132").Value = 14089
$ws.Range("I132").Value = 29223
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 87669
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -85139
$ws.Range("N132").Value = -17059.0001
# Row 134
$ws.Range("H134").Value = 2270.2856
$ws.Range("I134").Value = 2414.8
$ws.Range("J134").Value = 1909
$ws.Range("K134").Value = 7244.400000000001
$ws.Range("L134").Value = 5727
$ws.Range("M134").Value = -4709.400000000001
$ws.Range("N134").Value = -10797
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 675.8929000000001
$ws.Range("J113").Value = 675.8929000000001
$ws.Range("L113").Value = 2027.6787
$ws.Range("N113").Value = -6367.6787
# Row 131
$ws.Range("H131").Value = 21740564
$ws.Range("J131").Value = 1552.5238
$ws.Range("L131").Value = 4657.5714
$ws.Range("N131").Value = -14737.5714
# Row 133
$ws.Range("H133").Value = 3519
$ws.Range("I133").Value = 1529
$ws.Range("J133").Value = 4713
$ws.Range("K133").Value = 4587
$ws.Range("L133").Value = 14139
$ws.Range("M133").Value = 473
$ws.Range("N133").Value = -24259

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 8262650
$ws.Range("I11").Value = 8507812
$ws.Range("J11").Value = 7282000
$ws.Range("K11").Value = 8507812
$ws.Range("L11").Value = 7282000
$ws.Range("M11").Value = -8507673
$ws.Range("N11").Value = -7282278
# Row 21
$ws.Range("H21").Value = 2503200
$ws.Range("I21").Value = 5000000
$ws.Range("J21").Value = 6400
$ws.Range("K21").Value = 5000000
$ws.Range("L21").Value = 6400
$ws.Range("M21").Value = -4999827
$ws.Range("N21").Value = -6746
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 25
$ws.Range("H25").Value = 3663.3333
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# Row 30
$ws.Range("H30").Value = 2503200
$ws.Range("I30").Value = 5000000
$ws.Range("J30").Value = 6400
$ws.Range("K30").Value = 5000000
$ws.Range("L30").Value = 6400
$ws.Range("M30").Value = -4999895
$ws.Range("N30").Value = -6610

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1638.44
$ws.Range("I68").Value = 1530.1364
$ws.Range("J68").Value = 2432.6667
$ws.Range("K68").Value = 1530.1364
$ws.Range("L68").Value = 2432.6667
$ws.Range("M68").Value = -781.1364000000001
$ws.Range("N68").Value = -3930.6667
# Row 71
$ws.Range("H71").Value = 1638.44
$ws.Range("I71").Value = 1530.1364
$ws.Range("J71").Value = 2432.6667
$ws.Range("K71").Value = 7650.682000000001
$ws.Range("L71").Value = 12163.3335
$ws.Range("M71").Value = -3906.682000000001
$ws.Range("N71").Value = -19651.3335

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1886.3636
$ws.Range("I132").Value = 1374.0358
$ws.Range("J132").Value = 4755.4
$ws.Range("K132").Value = 4122.107400000001
$ws.Range("L132").Value = 14266.2
$ws.Range("M132").Value = -1592.107400000001
$ws.Range("N132").Value = -19326.2
